$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose values look numeric, so Excel
# preserves them as literal text (matching the source feeds display strings)
# instead of converting to floating-point numbers.
$textCells = @("D5","D6","D9","D10","D11","D12","D13","D17","D18","D20","D21","D23","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D42","D43","D44","D47","D48","D49","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "73.112.81"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").Value = "3.988.10"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "616.46"
$ws.Range("E5").Value = "  +14.71%  "
$ws.Range("D6").Value = "166.43"
$ws.Range("E6").Value = "  +11.39%  "
$ws.Range("E7").Value = "  -1.33%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "0.758"
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("D11").Value = "57.33"
$ws.Range("E11").Value = "  +6.78%  "
$ws.Range("D12").Value = "0.0000331"
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").Value = "11.14"
$ws.Range("E13").Value = "  +2.67%  "
$ws.Range("D14").Value = "4.621.04"
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("D15").Value = "3.991.90"
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("E16").Value = "  +5.27%  "
$ws.Range("D17").Value = "14.25"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "20.58"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "72.992.43"
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").Value = "0.131"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").Value = "439.73"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("E22").Value = "  +16.14%  "
$ws.Range("D23").Value = "96.26"
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("E24").Value = "  -4.25%  "
$ws.Range("D25").Value = "14.30"
$ws.Range("E25").Value = "  -2.52%  "
$ws.Range("D26").Value = "4.10"
$ws.Range("E26").Value = "  -7.54%  "
$ws.Range("D27").Value = "11.17"
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "5.96"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "10.53"
$ws.Range("E29").Value = "  -2.06%  "
$ws.Range("D30").Value = "36.21"
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("D31").Value = "7.79"
$ws.Range("E31").Value = "  -5.53%  "
$ws.Range("D32").Value = "13.79"
$ws.Range("E32").Value = "  +1.24%  "
$ws.Range("D33").Value = "0.130"
$ws.Range("E33").Value = "  -3.83%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "72.45"
$ws.Range("E34").Value = "  +7.30%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "48.63"
$ws.Range("E35").Value = "  -3.37%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0000102"
$ws.Range("E36").Value = "  +23.65%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "639.49"
$ws.Range("E37").Value = "  -6.35%  "
$ws.Range("D38").Value = "0.433"
$ws.Range("E38").Value = "  -5.92%  "
$ws.Range("D39").Value = "3.48"
$ws.Range("E39").Value = "  +2.99%  "
$ws.Range("E40").Value = "  -1.34%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "11.10"
$ws.Range("E42").Value = "  -1.10%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "3.30"
$ws.Range("E43").Value = "  -3.50%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("E45").Value = "  -1.28%  "
$ws.Range("E46").Value = "  -1.00%  "
$ws.Range("D47").Value = "3.41"
$ws.Range("E47").Value = "  +1.76%  "
$ws.Range("D48").Value = "2.61"
$ws.Range("E48").Value = "  -2.21%  "
$ws.Range("D49").Value = "2.87"
$ws.Range("E49").Value = "  +31.78%  "
$ws.Range("D50").Value = "2.876.48"
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("D51").Value = "3.04"
$ws.Range("E51").Value = "  -2.01%  "
